$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to reflect custom-accuracy re-export (rounded values)
$ws.Cells.Item(5, 1).Value = 44781.9027662037
$ws.Cells.Item(5, 2).Value = 4.8
$ws.Cells.Item(5, 3).Value = 3.27
$ws.Cells.Item(5, 4).Value = 0.72
$ws.Cells.Item(5, 5).Value = 10.14
$ws.Cells.Item(5, 6).Value = 8.390000000000001
$ws.Cells.Item(5, 7).Value = 3.78
$ws.Cells.Item(5, 8).Value = 19.19
$ws.Cells.Item(5, 9).Value = 5.82
$ws.Cells.Item(5, 10).Value = 2.5
$ws.Cells.Item(5, 11).Value = 3.7
$ws.Cells.Item(5, 12).Value = 4.16
$ws.Cells.Item(5, 13).Value = 4.23
$ws.Cells.Item(5, 14).Value = 1.21
$ws.Cells.Item(5, 15).Value = 3.76
$ws.Cells.Item(5, 16).Value = 5.31
$ws.Cells.Item(5, 17).Value = 3.35
$ws.Cells.Item(5, 18).Value = 0.73
$ws.Cells.Item(5, 19).Value = 0.37
$ws.Cells.Item(5, 20).Value = 49.94
$ws.Cells.Item(5, 21).Value = 10.74
$ws.Cells.Item(5, 22).Value = 3.47
$ws.Cells.Item(5, 23).Value = 7.05
$ws.Cells.Item(5, 24).Value = 3.84
$ws.Cells.Item(5, 25).Value = 0.39
$ws.Cells.Item(5, 26).Value = 8.69
$ws.Cells.Item(5, 27).Value = 3.07
$ws.Cells.Item(5, 28).Value = 2.84
$ws.Cells.Item(5, 29).Value = 3.31
$ws.Cells.Item(5, 30).Value = 4.3
$ws.Cells.Item(5, 31).Value = 0.5600000000000001
$ws.Cells.Item(5, 32).Value = 17.53
$ws.Cells.Item(5, 33).Value = 1.87
$ws.Cells.Item(5, 34).Value = 4.34

# Remove row 6 entirely (dataset trimmed)
$ws.Rows.Item(6).Delete() | Out-Null

